$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 464.25
$ws.Range("J17").Value = 451.85184
$ws.Range("L17").Value = 1355.55552
$ws.Range("N17").Value = -1691.55552
$ws.Range("H19").Value = 991.4286
$ws.Range("J19").Value = 1060.75
$ws.Range("L19").Value = 1060.75
$ws.Range("N19").Value = -1410.75
$ws.Range("H33").Value = 200.70589
$ws.Range("I33").Value = 200.70589
$ws.Range("K33").Value = 200.70589
$ws.Range("M33").Value = 28.29410999999999
$ws.Range("H88").Value = 5557607
$ws.Range("I88").Value = 11112450
$ws.Range("J88").Value = 2763.3333
$ws.Range("K88").Value = 11112450
$ws.Range("L88").Value = 2763.3333
$ws.Range("M88").Value = -11112044
$ws.Range("N88").Value = -3575.3333
$ws.Range("H91").Value = 5557607
$ws.Range("I91").Value = 11112450
$ws.Range("J91").Value = 2763.3333
$ws.Range("K91").Value = 11112450
$ws.Range("L91").Value = 2763.3333
$ws.Range("M91").Value = -11111046
$ws.Range("N91").Value = -5571.3333
$ws.Range("H92").Value = 1027.8
$ws.Range("I92").Value = 1074.75
$ws.Range("J92").Value = 840
$ws.Range("K92").Value = 1074.75
$ws.Range("L92").Value = 840
$ws.Range("M92").Value = 173.25
$ws.Range("N92").Value = -3336
$ws.Range("H98").Value = 3632.2104
$ws.Range("I98").Value = 4025.0625
$ws.Range("J98").Value = 1537
$ws.Range("K98").Value = 4025.0625
$ws.Range("L98").Value = 1537
$ws.Range("M98").Value = -2527.0625
$ws.Range("N98").Value = -4533
$ws.Range("H116").Value = 8532
$ws.Range("I116").Value = 4048
$ws.Range("K116").Value = 4048
$ws.Range("M116").Value = -606
$ws.Range("H122").Value = 3632.2104
$ws.Range("I122").Value = 4025.0625
$ws.Range("J122").Value = 1537
$ws.Range("K122").Value = 12075.1875
$ws.Range("L122").Value = 4611
$ws.Range("M122").Value = -9625.1875
$ws.Range("N122").Value = -9511
$ws.Range("H136").Value = 84366
$ws.Range("J136").Value = 83694.5
$ws.Range("L136").Value = 83694.5
$ws.Range("N136").Value = -93894.5
$ws.Range("H137").Value = 1934.1111
$ws.Range("I137").Value = 1973.1428
$ws.Range("K137").Value = 5919.428400000001
$ws.Range("M137").Value = -3369.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1540431.6
$ws.Range("I32").Value = 727484.5600000001
$ws.Range("K32").Value = 727484.5600000001
$ws.Range("M32").Value = -727197.5600000001
$ws.Range("H61").Value = 3854.5557
$ws.Range("I61").Value = 3346
$ws.Range("K61").Value = 3346
$ws.Range("M61").Value = -3134
$ws.Range("H135").Value = 93781.164
$ws.Range("J135").Value = 93781.164
$ws.Range("L135").Value = 93781.164
$ws.Range("N135").Value = -103921.164
$ws.Range("H136").Value = 3854.5557
$ws.Range("I136").Value = 3346
$ws.Range("K136").Value = 10038
$ws.Range("M136").Value = -7488

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3124.027
$ws.Range("I86").Value = 3021.1333
$ws.Range("J86").Value = 3565
$ws.Range("K86").Value = 3021.1333
$ws.Range("L86").Value = 3565
$ws.Range("M86").Value = -1898.1333
$ws.Range("N86").Value = -5811
$ws.Range("H89").Value = 3124.027
$ws.Range("I89").Value = 3021.1333
$ws.Range("J89").Value = 3565
$ws.Range("K89").Value = 15105.6665
$ws.Range("L89").Value = 17825
$ws.Range("M89").Value = -9489.666499999999
$ws.Range("N89").Value = -29057
$ws.Range("H94").Value = 222234050
$ws.Range("I94").Value = 285728930
$ws.Range("K94").Value = 285728930
$ws.Range("M94").Value = -285728479
$ws.Range("H134").Value = 1790.2142
$ws.Range("I134").Value = 1131.35
$ws.Range("K134").Value = 3394.05
$ws.Range("M134").Value = -859.0499999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5437680.5
$ws.Range("I31").Value = 2049.5293
$ws.Range("K31").Value = 2049.5293
$ws.Range("M31").Value = -1754.5293
$ws.Range("H34").Value = 5437680.5
$ws.Range("I34").Value = 2049.5293
$ws.Range("K34").Value = 2049.5293
$ws.Range("M34").Value = -1847.5293
$ws.Range("H99").Value = 4176.385
$ws.Range("I99").Value = 3517.2856
$ws.Range("J99").Value = 4945.3335
$ws.Range("K99").Value = 3517.2856
$ws.Range("L99").Value = 4945.3335
$ws.Range("M99").Value = -2019.2856
$ws.Range("N99").Value = -7941.3335
$ws.Range("H122").Value = 2148.7917
$ws.Range("J122").Value = 4075.5
$ws.Range("L122").Value = 12226.5
$ws.Range("N122").Value = -17126.5
$ws.Range("H126").Value = 4176.385
$ws.Range("I126").Value = 3517.2856
$ws.Range("J126").Value = 4945.3335
$ws.Range("K126").Value = 10551.8568
$ws.Range("L126").Value = 14836.0005
$ws.Range("M126").Value = -8081.856800000001
$ws.Range("N126").Value = -19776.0005
$ws.Range("H132").Value = 3824.1943
$ws.Range("I132").Value = 3226.88
$ws.Range("K132").Value = 9680.639999999999
$ws.Range("M132").Value = -7150.639999999999
$ws.Range("H134").Value = 3457.9756
$ws.Range("J134").Value = 3389.7778
$ws.Range("L134").Value = 10169.3334
$ws.Range("N134").Value = -15239.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11251035
$ws.Range("I4").Value = 12466141
$ws.Range("J4").Value = 5175505
$ws.Range("K4").Value = 37398423
$ws.Range("L4").Value = 15526515
$ws.Range("M4").Value = -37398311
$ws.Range("N4").Value = -15526739
$ws.Range("H12").Value = 39.4
$ws.Range("J12").Value = 49
$ws.Range("L12").Value = 147
$ws.Range("N12").Value = -493
$ws.Range("H113").Value = 1117.6666
$ws.Range("I113").Value = 889.6667
$ws.Range("J113").Value = 1231.6666
$ws.Range("K113").Value = 2669.0001
$ws.Range("L113").Value = 3694.9998
$ws.Range("M113").Value = -499.0001000000002
$ws.Range("N113").Value = -8034.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 90911740
$ws.Range("I80").Value = 200002450
$ws.Range("J80").Value = 2816.6667
$ws.Range("K80").Value = 200002450
$ws.Range("L80").Value = 2816.6667
$ws.Range("M80").Value = -200001452
$ws.Range("N80").Value = -4812.6667
$ws.Range("H83").Value = 90911740
$ws.Range("I83").Value = 200002450
$ws.Range("J83").Value = 2816.6667
$ws.Range("K83").Value = 1000012250
$ws.Range("L83").Value = 14083.3335
$ws.Range("M83").Value = -1000007258
$ws.Range("N83").Value = -24067.3335
$ws.Range("H122").Value = 3037.8667
$ws.Range("I122").Value = 2649.4211
$ws.Range("J122").Value = 3708.818
$ws.Range("K122").Value = 7948.263300000001
$ws.Range("L122").Value = 11126.454
$ws.Range("M122").Value = -5498.263300000001
$ws.Range("N122").Value = -16026.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 22969.47
$ws.Range("I40").Value = 24092.562
$ws.Range("K40").Value = 24092.562
$ws.Range("M40").Value = -23956.562
$ws.Range("H82").Value = 577.7857
$ws.Range("I82").Value = 544.5454999999999
$ws.Range("K82").Value = 544.5454999999999
$ws.Range("M82").Value = -183.5454999999999
$ws.Range("H85").Value = 577.7857
$ws.Range("I85").Value = 544.5454999999999
$ws.Range("K85").Value = 544.5454999999999
$ws.Range("M85").Value = 703.4545000000001
$ws.Range("H93").Value = 2430.2
$ws.Range("I93").Value = 2143.1482
$ws.Range("J93").Value = 3399
$ws.Range("K93").Value = 2143.1482
$ws.Range("L93").Value = 3399
$ws.Range("M93").Value = -895.1482000000001
$ws.Range("N93").Value = -5895

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 545
$ws.Range("I107").Value = 533.375
$ws.Range("J107").Value = 591.5
$ws.Range("K107").Value = 1600.125
$ws.Range("L107").Value = 1774.5
$ws.Range("M107").Value = 319.875
$ws.Range("N107").Value = -5614.5
$ws.Range("H132").Value = 2524.8
$ws.Range("I132").Value = 2592.5417
$ws.Range("K132").Value = 7777.625100000001
$ws.Range("M132").Value = -5247.625100000001
$ws.Range("H136").Value = 3366.2942
$ws.Range("I136").Value = 3569.1333
$ws.Range("J136").Value = 1845
$ws.Range("K136").Value = 10707.3999
$ws.Range("L136").Value = 5535
$ws.Range("M136").Value = -8157.3999
$ws.Range("N136").Value = -10635
